$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: reset to unknown/placeholder values
$ws.Range("C2").Value = "Unknown Title"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "1970-01-01"

# Row 3: reset to unknown/placeholder values
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "1970-01-01"

# Row 4: reset to unknown/placeholder values
$ws.Range("C4").Value = "Unknown Title"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "not found"
$ws.Range("G4").Value = "N/A"
$ws.Range("I4").Value = ""

# Row 5: reset to unknown/placeholder values
$ws.Range("C5").Value = "Unknown Title"
$ws.Range("E5").Value = "[]"
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"
$ws.Range("H5").Value = "1970-01-01"

# Row 6: reset to unknown/placeholder values
$ws.Range("C6").Value = "Unknown Title"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"
$ws.Range("H6").Value = "1970-01-01"

# Row 7: reset to unknown/placeholder values
$ws.Range("C7").Value = "Unknown Title"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = "not found"
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = "1970-01-01"

# Row 8: reset to unknown/placeholder values
$ws.Range("C8").Value = "Unknown Title"
$ws.Range("D8").Value = "Unknown Abstract"
$ws.Range("E8").Value = "[]"
$ws.Range("F8").Value = "not found"
$ws.Range("G8").Value = "N/A"
$ws.Range("H8").Value = "1970-01-01"
$ws.Range("I8").Value = ""

# Row 9: reset to unknown/placeholder values
$ws.Range("C9").Value = "Unknown Title"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "not found"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = "1970-01-01"

# Row 10: reset to unknown/placeholder values
$ws.Range("C10").Value = "Unknown Title"
$ws.Range("E10").Value = "[]"
$ws.Range("F10").Value = "not found"
$ws.Range("G10").Value = "N/A"
$ws.Range("H10").Value = "1970-01-01"

# Row 14: reset to unknown/placeholder values
$ws.Range("C14").Value = "Unknown Title"
$ws.Range("E14").Value = "[]"
$ws.Range("F14").Value = "not found"
$ws.Range("G14").Value = "N/A"
$ws.Range("H14").Value = "1970-01-01"

# Row 16: reset to unknown/placeholder values
$ws.Range("C16").Value = "Unknown Title"
$ws.Range("E16").Value = "[]"
$ws.Range("F16").Value = "not found"
$ws.Range("G16").Value = "N/A"
$ws.Range("H16").Value = "1970-01-01"

# Row 17: reset to unknown/placeholder values
$ws.Range("C17").Value = "Unknown Title"
$ws.Range("E17").Value = "[]"
$ws.Range("F17").Value = "not found"
$ws.Range("G17").Value = "N/A"
$ws.Range("H17").Value = "1970-01-01"

# Row 18: reset to unknown/placeholder values
$ws.Range("C18").Value = "Unknown Title"
$ws.Range("D18").Value = "Unknown Abstract"
$ws.Range("E18").Value = "[]"
$ws.Range("F18").Value = "not found"
$ws.Range("G18").Value = "N/A"
$ws.Range("H18").Value = "1970-01-01"
$ws.Range("I18").Value = ""

# Row 11: update Authors and Misc. Data with full text
$ws.Range("E11").Value = "[Joelle%O’Neil%NULL%0, Susan%Hughes%susan.hughes@fresno.ucsf.edu%1, Andrea%Lourie%NULL%1, John%Zweifler%NULL%1]"
$ws.Range("J11").Value = "American College of Allergy, Asthma &amp; Immunology. Published by Elsevier Inc."

# Row 12: update Abstract, Authors, and Misc. Data with full text
$ws.Range("D12").Value = "`nObjective.`n To investigate the safety (risk) and efficacy (benefit) of Echinacea purpurea extract in the prevention of common cold episodes in a large population over a 4-month period.`n Methods.`n 755 healthy subjects were allocated to receive either an alcohol extract from freshly harvested E.`n purpurea (95% herba and 5% root) or placebo.`n Participants were required to record adverse events and to rate cold-related issues in a diary throughout the investigation period.`n Nasal secretions were sampled at acute colds and screened for viruses.`n Results.`n A total of 293 adverse events occurred with Echinacea and 306 with placebo treatment.`n Nine and 10% of participants experienced adverse events, which were at least possibly related to the study drug (adverse drug reactions).`n Thus, the safety of Echinacea was noninferior to placebo.`n Echinacea reduced the total number of cold episodes, cumulated episode days within the group, and pain-killer medicated episodes.`n Echinacea inhibited virally confirmed colds and especially prevented enveloped virus infections (P &lt; 0.05).`n Echinacea showed maximal effects on recurrent infections, and preventive effects increased with therapy compliance and adherence to the protocol.`n Conclusions.`n Compliant prophylactic intake of E.`n purpurea over a 4-month period appeared to provide a positive risk to benefit ratio.`n"
$ws.Range("E12").Value = "[M.%Jawad%NULL%0, R.%Schoop%NULL%1, A.%Suter%NULL%1, P.%Klein%NULL%1, R.%Eccles%NULL%1]"
$ws.Range("J12").Value = "Hindawi Publishing Corporation"

# Row 13: update Abstract, Authors, and Misc. Data with full text
$ws.Range("D13").Value = "`nObjective.`n To identify whether a standardised Echinacea formulation is effective in the prevention of respiratory and other symptoms associated with long-haul flights.`n Methods.`n 175 adults participated in a randomised, double-blind placebo-controlled trial travelling back from Australia to America, Europe, or Africa for a period of 1–5 weeks on commercial flights via economy class.`n Participants took Echinacea (root extract, standardised to 4.4 mg alkylamides) or placebo tablets.`n Participants were surveyed before, immediately after travel, and at 4 weeks after travel regarding upper respiratory symptoms and travel-related quality of life.`n Results.`n Respiratory symptoms for both groups increased significantly during travel (P &lt; 0.0005).`n However, the Echinacea group had borderline significantly lower respiratory symptom scores compared to placebo (P = 0.05) during travel.`n Conclusions.`n Supplementation with standardised Echinacea tablets, if taken before and during travel, may have preventive effects against the development of respiratory symptoms during travel involving long-haul flights.`n"
$ws.Range("E13").Value = "[E.%Tiralongo%NULL%0, R. A.%Lea%NULL%1, S. S.%Wee%NULL%1, M. M.%Hanna%NULL%1, L. R.%Griffiths%NULL%1]"
$ws.Range("J13").Value = "Hindawi Publishing Corporation"

# Row 15: update Abstract, Authors, and Misc. Data with full text
$ws.Range("D15").Value = "A randomized, double-blind, placebo-controlled clinical trial was conducted to evaluate the ability of Echinacea purpurea to prevent infection with rhinovirus type 39 (RV-39).`n Forty-eight previously healthy adults received echinacea or placebo, 2.5 mL 3 times per day, for 7 days before and 7 days after intranasal inoculation with RV-39. Symptoms were assessed to evaluate clinical illness.`n Viral culture and serologic studies were performed to evaluate the presence of rhinovirus infection.`n A total of 92% of echinacea recipients and 95% of placebo recipients were infected.`n Colds developed in 58% of echinacea recipients and 82% of placebo recipients (P = .`n114, by Fisher's exact test).`n Administration of echinacea before and after exposure to rhinovirus did not decrease the rate of infection; however, because of the small sample size, statistical hypothesis testing had relatively poor power to detect statistically significant differences in the frequency and severity of illness.`n"
$ws.Range("E15").Value = "[Steven J.%Sperber%ssperber@humed.com%0, Leena P.%Shah%NULL%1, Richard D.%Gilbert%NULL%1, Thomas W.%Ritchey%NULL%1, Arnold S.%Monto%NULL%1]"
$ws.Range("J15").Value = "The University of Chicago Press"
